$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "public budget" data replaces the old "agricultural machinery" data.
# Columns: A=input, B=chn_block4, C=asis, D=variables
$data = @(
    @("合计",     "v6_cz_yszc_hj"),
    @("教育",     "v6_cz_yszc_jy"),
    @("科学技术", "v6_cz_yszc_kxjs"),
    @("农林水",   "v6_cz_yszc_nls")
)

# Clear out the old rows (2..12) first so no stale rows/cells remain.
$ws.Range("A2:D12").Clear()

$r = 2
foreach ($row in $data) {
    $chn = $row[0]
    $var = $row[1]
    $ws.Cells.Item($r, 1).Value = $chn
    $ws.Cells.Item($r, 2).Value = $chn
    $ws.Cells.Item($r, 3).Value = $true
    $ws.Cells.Item($r, 4).Value = $var
    $r++
}
